$wb = $excel.ActiveWorkbook
$rawWs = $wb.Worksheets.Item("Raw Data")
$tableWs = $wb.Worksheets.Item("Tables With Automatic Formulas")

# --- Raw Data sheet: new column N "VaccinationRate" ---
$rawWs.Range("N1").Value = "VaccinationRate"

$rawWs.Range("N2").Value = 0.31273322999999997
$rawWs.Range("N3").Value = 1.6062541530000001
$rawWs.Range("N4").Value = 0.32974390999999997
$rawWs.Range("N5").Value = 0.53762031600000004
$rawWs.Range("N6").Value = 0.32410455700000002
$rawWs.Range("N7").Value = 1.1102390680000001
$rawWs.Range("N8").Value = 0.46148747499999998
$rawWs.Range("N9").Value = 0.74445535900000004
$rawWs.Range("N10").Value = 0.64056939499999999
$rawWs.Range("N11").Value = 0.55759579000000004
$rawWs.Range("N12").Value = 0.24725594400000001
$rawWs.Range("N14").Value = 0.57569089699999998
$rawWs.Range("N15").Value = 0.80211353900000004
$rawWs.Range("N16").Value = 0.59773309699999999
$rawWs.Range("N17").Value = 0.26615269600000002
$rawWs.Range("N18").Value = 0.063129512999999998
$rawWs.Range("N19").Value = 0.250469518
$rawWs.Range("N20").Value = 0.70314813600000003
$rawWs.Range("N21").Value = 1.275203138
$rawWs.Range("N22").Value = 0.327961634
$rawWs.Range("N23").Value = 0.51603268300000005
$rawWs.Range("N24").Value = 0.37674350200000001
$rawWs.Range("N25").Value = 0.20633291300000001
$rawWs.Range("N26").Value = 0.150679111
$rawWs.Range("N27").Value = 0.37541716200000003
$rawWs.Range("N28").Value = 0.69724425199999995
$rawWs.Range("N29").Value = 0.80145840300000004
$rawWs.Range("N30").Value = 0.50115657899999999
$rawWs.Range("N31").Value = 0.28153750900000002
$rawWs.Range("N32").Value = 0.40971121100000002
$rawWs.Range("N33").Value = 0.66811743599999995
$rawWs.Range("N34").Value = 0.58058538999999998
$rawWs.Range("N35").Value = 0.235954155
$rawWs.Range("N36").Value = 1.7090666050000001
$rawWs.Range("N37").Value = 0.44241984699999998
$rawWs.Range("N38").Value = 0.75385977500000001
$rawWs.Range("N39").Value = 0.34657661299999998
$rawWs.Range("N40").Value = 0.44975648800000001
$rawWs.Range("N41").Value = 0.66886089400000004
$rawWs.Range("N42").Value = 0.38638829000000002
$rawWs.Range("N43").Value = 1.554020049
$rawWs.Range("N44").Value = 0.78668929600000004
$rawWs.Range("N45").Value = 0.48689645500000001
$rawWs.Range("N46").Value = 0.53752828100000005
$rawWs.Range("N47").Value = 0.69838623899999996
$rawWs.Range("N48").Value = 0.448155146
$rawWs.Range("N49").Value = 0.39810945800000003
$rawWs.Range("N50").Value = 1.643312342
$rawWs.Range("N51").Value = 0.17816610399999999
$rawWs.Range("N52").Value = 0.49589609200000001

$rawWs.Range("N13").Value = "–"

$excel.Calculate()

# --- Tables sheet: new column H "Vaccination Rate (%)" ---
$tableWs.Range("G1").Copy()
$tableWs.Range("H1").PasteSpecial(-4122)
$tableWs.Range("H1").Value = "Vaccination Rate (%)"
$tableWs.Range("H1").Borders.Item(8).LineStyle = -4142
$tableWs.Range("H1").Borders.Item(9).LineStyle = -4142

$tableWs.Range("H2").Formula = "=VLOOKUP(A2,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H3").Formula = "=VLOOKUP(A3,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H4").Formula = "=VLOOKUP(A4,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H5").Formula = "=VLOOKUP(A5,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H6").Formula = "=VLOOKUP(A6,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H7").Formula = "=VLOOKUP(A7,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H8").Formula = "=VLOOKUP(A8,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H9").Formula = "=VLOOKUP(A9,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H10").Formula = "=VLOOKUP(A10,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H11").Formula = "=VLOOKUP(A11,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H12").Formula = "=VLOOKUP(A12,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H13").Formula = "=VLOOKUP(A13,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H14").Formula = "=VLOOKUP(A14,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H15").Formula = "=VLOOKUP(A15,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H16").Formula = "=VLOOKUP(A16,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H17").Formula = "=VLOOKUP(A17,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H18").Formula = "=VLOOKUP(A18,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H19").Formula = "=VLOOKUP(A19,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H20").Formula = "=VLOOKUP(A20,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H21").Formula = "=VLOOKUP(A21,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H22").Formula = "=VLOOKUP(A22,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H23").Formula = "=VLOOKUP(A23,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H24").Formula = "=VLOOKUP(A24,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H25").Formula = "=VLOOKUP(A25,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H26").Formula = "=VLOOKUP(A26,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H27").Formula = "=VLOOKUP(A27,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H28").Formula = "=VLOOKUP(A28,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H29").Formula = "=VLOOKUP(A29,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H30").Formula = "=VLOOKUP(A30,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H31").Formula = "=VLOOKUP(A31,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H32").Formula = "=VLOOKUP(A32,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H33").Formula = "=VLOOKUP(A33,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H34").Formula = "=VLOOKUP(A34,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H35").Formula = "=VLOOKUP(A35,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H36").Formula = "=VLOOKUP(A36,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H37").Formula = "=VLOOKUP(A37,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H38").Formula = "=VLOOKUP(A38,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H39").Formula = "=VLOOKUP(A39,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H40").Formula = "=VLOOKUP(A40,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H41").Formula = "=VLOOKUP(A41,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H42").Formula = "=VLOOKUP(A42,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H43").Formula = "=VLOOKUP(A43,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H44").Formula = "=VLOOKUP(A44,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H45").Formula = "=VLOOKUP(A45,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H46").Formula = "=VLOOKUP(A46,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H47").Formula = "=VLOOKUP(A47,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H48").Formula = "=VLOOKUP(A48,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H49").Formula = "=VLOOKUP(A49,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H50").Formula = "=VLOOKUP(A50,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H51").Formula = "=VLOOKUP(A51,'Raw Data'!A:N,14,FALSE)"
$tableWs.Range("H52").Formula = "=VLOOKUP(A52,'Raw Data'!A:N,14,FALSE)"

$excel.Calculate()

# --- Restore selections to match final authored state ---
$rawWs.Activate()
$rawWs.Range("E16").Select()
$tableWs.Activate()
$tableWs.Range("F5").Select()
